# The workbook is already open; grab the active workbook/sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Update tab names in all BOMs" -> rename this sheet's tab
# from "Nucleus" to "BOM".
$ws.Name = "BOM"

# Scroll the view so row 3 is the first visible row (topLeftCell = A3),
# then move the active cell/selection down to A47 (where the user was
# last working), matching the updated cursor position in the saved file.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A47").Select() | Out-Null
